$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1 - copy formatting from H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Data values for columns I and J, rows 2-32
$iValues = @(6,7,8,8,9,7,8,6,9,6,7,6,9,6,8,7,8,8,5,8,8,7,7,9,6,6,4,6,6,4,4)
$jValues = @(8,8,8,9,9,7,9,7,9,7,8,8,9,8,8,7,8,9,5,8,8,7,7,9,7,7,4,7,6,4,4)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
